$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "yes" values in column C for rows 5-10 (the "ignore" column),
# removing the cells entirely so they no longer appear in the sheet data
# and the now-unused "yes" shared string is dropped.
$ws.Range("C5:C10").ClearContents()
